$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 88846.664
$ws.Range("I62").Value = 170067.33
$ws.Range("J62").Value = 7626
$ws.Range("K62").Value = 170067.33
$ws.Range("L62").Value = 7626
$ws.Range("M62").Value = -169443.33
$ws.Range("N62").Value = -8874
$ws.Range("H65").Value = 88846.664
$ws.Range("I65").Value = 170067.33
$ws.Range("J65").Value = 7626
$ws.Range("K65").Value = 850336.6499999999
$ws.Range("L65").Value = 38130
$ws.Range("M65").Value = -847216.6499999999
$ws.Range("N65").Value = -44370
$ws.Range("H112").Value = 1380.7727
$ws.Range("J112").Value = 1435.1052
$ws.Range("L112").Value = 4305.3156
$ws.Range("N112").Value = -6521.3156
$ws.Range("H116").Value = 4417.857
$ws.Range("I116").Value = 4372.222
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 4372.222
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = -930.2219999999998
$ws.Range("N116").Value = -11384
$ws.Range("H121").Value = 859
$ws.Range("I121").Value = 565
$ws.Range("J121").Value = 1300
$ws.Range("K121").Value = 1695
$ws.Range("L121").Value = 3900
$ws.Range("M121").Value = 52
$ws.Range("N121").Value = -7394
$ws.Range("H129").Value = 1425420.2
$ws.Range("I129").Value = 281
$ws.Range("J129").Value = 2470522.5
$ws.Range("K129").Value = 843
$ws.Range("L129").Value = 7411567.5
$ws.Range("M129").Value = 4157
$ws.Range("N129").Value = -7421567.5
$ws.Range("H131").Value = 1609.2941
$ws.Range("I131").Value = 670.25
$ws.Range("J131").Value = 3863
$ws.Range("K131").Value = 2010.75
$ws.Range("L131").Value = 11589
$ws.Range("M131").Value = 3029.25
$ws.Range("N131").Value = -21669
$ws.Range("H132").Value = 2782.111
$ws.Range("I132").Value = 2941.1365
$ws.Range("K132").Value = 8823.4095
$ws.Range("M132").Value = -6293.4095
$ws.Range("H137").Value = 798.9091
$ws.Range("I137").Value = 798.9091
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2396.7273
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 153.2727
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 3479.23
$ws.Range("I138").Value = 1555.6428
$ws.Range("J138").Value = 3792.372
$ws.Range("K138").Value = 4666.928400000001
$ws.Range("L138").Value = 11377.116
$ws.Range("M138").Value = 473.0715999999993
$ws.Range("N138").Value = -21657.116

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5822.5
$ws.Range("I32").Value = 4113.0684
$ws.Range("J32").Value = 18358.334
$ws.Range("K32").Value = 4113.0684
$ws.Range("L32").Value = 18358.334
$ws.Range("M32").Value = -3826.0684
$ws.Range("N32").Value = -18932.334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2779.5557
$ws.Range("I31").Value = 2846.5435
$ws.Range("J31").Value = 2394.375
$ws.Range("K31").Value = 2846.5435
$ws.Range("L31").Value = 2394.375
$ws.Range("M31").Value = -2551.5435
$ws.Range("N31").Value = -2984.375
$ws.Range("H34").Value = 2779.5557
$ws.Range("I34").Value = 2846.5435
$ws.Range("J34").Value = 2394.375
$ws.Range("K34").Value = 2846.5435
$ws.Range("L34").Value = 2394.375
$ws.Range("M34").Value = -2644.5435
$ws.Range("N34").Value = -2798.375
$ws.Range("H52").Value = 26300
$ws.Range("J52").Value = 26300
$ws.Range("L52").Value = 26300
$ws.Range("N52").Value = -26888
$ws.Range("H86").Value = 12239.8
$ws.Range("I86").Value = 12258
$ws.Range("J86").Value = 12221.6
$ws.Range("K86").Value = 12258
$ws.Range("L86").Value = 12221.6
$ws.Range("M86").Value = -11135
$ws.Range("N86").Value = -14467.6
$ws.Range("H89").Value = 12239.8
$ws.Range("I89").Value = 12258
$ws.Range("J89").Value = 12221.6
$ws.Range("K89").Value = 61290
$ws.Range("L89").Value = 61108
$ws.Range("M89").Value = -55674
$ws.Range("N89").Value = -72340
$ws.Range("H103").Value = 12333.333
$ws.Range("I103").Value = 3500
$ws.Range("K103").Value = 3500
$ws.Range("M103").Value = -2328
$ws.Range("H132").Value = 1516
$ws.Range("I132").Value = 1144.4615
$ws.Range("J132").Value = 2723.5
$ws.Range("K132").Value = 3433.3845
$ws.Range("L132").Value = 8170.5
$ws.Range("M132").Value = -903.3844999999997
$ws.Range("N132").Value = -13230.5
$ws.Range("H134").Value = 16130289
$ws.Range("I134").Value = 1125.3684
$ws.Range("J134").Value = 41668132
$ws.Range("K134").Value = 3376.1052
$ws.Range("L134").Value = 125004396
$ws.Range("M134").Value = -841.1052
$ws.Range("N134").Value = -125009466

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 840.4
$ws.Range("J5").Value = 845
$ws.Range("L5").Value = 2535
$ws.Range("N5").Value = -2759
$ws.Range("H56").Value = 3364.3333
$ws.Range("I56").Value = 3364.3333
$ws.Range("K56").Value = 3364.3333
$ws.Range("M56").Value = -2834.3333
$ws.Range("H105").Value = 154957.14
$ws.Range("J105").Value = 154957.14
$ws.Range("L105").Value = 464871.42
$ws.Range("N105").Value = -470113.42
$ws.Range("H112").Value = 2676.7273
$ws.Range("I112").Value = 1157.3334
$ws.Range("K112").Value = 3472.0002
$ws.Range("M112").Value = -2364.0002
$ws.Range("H135").Value = 840.4
$ws.Range("J135").Value = 845
$ws.Range("L135").Value = 7605
$ws.Range("N135").Value = -12675

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6233.3335
$ws.Range("J70").Value = 5675
$ws.Range("L70").Value = 5675
$ws.Range("N70").Value = -6215
$ws.Range("H73").Value = 6233.3335
$ws.Range("J73").Value = 5675
$ws.Range("L73").Value = 5675
$ws.Range("N73").Value = -7547
$ws.Range("H98").Value = 27920
$ws.Range("J98").Value = 27920
$ws.Range("L98").Value = 27920
$ws.Range("N98").Value = -33910
$ws.Range("H132").Value = 2589.25
$ws.Range("I132").Value = 2306.6428
$ws.Range("K132").Value = 6919.928400000001
$ws.Range("M132").Value = -4389.928400000001
$ws.Range("H134").Value = 18645.2
$ws.Range("J134").Value = 18645.2
$ws.Range("L134").Value = 55935.60000000001
$ws.Range("N134").Value = -61005.60000000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -9480
$ws.Range("N54").ClearContents()
$ws.Range("H113").Value = 597.15
$ws.Range("I113").Value = 459.16666
$ws.Range("K113").Value = 1377.49998
$ws.Range("M113").Value = 792.5000199999999
